# Update WG attendees minutes with the latest working-group meeting:
# add Margaret Wishart (Bristol Myers Squibb) as a new attendee row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new attendee row (row 18) ---------------------------------
$ws.Range("A18").Value = "Margaret Wishart"
$ws.Range("B18").Value = "Bristol Myers Squibb"
$ws.Range("C18").Value = "https://www.bms.com/ "

# Match the formatting used for the other Bristol Myers Squibb row (row 5)
# so the new C18 cell gets the hyperlink-blue style instead of the plain
# default style that a freshly written cell otherwise inherits.
$ws.Range("C5").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# --- Hyperlink the affiliation URL cell ---------------------------------
$ws.Hyperlinks.Add($ws.Range("C18"), "https://www.bms.com/", "", "", "https://www.bms.com/ ")

# --- Re-style the two rows that used to carry the (duplicate) blue link
#     font so they point at the same shared hyperlink style as the rest
#     of the affiliation_url column. -------------------------------------
$ws.Range("C2").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)

# --- Move the active selection down to the new empty row, like a user
#     would leave it after typing the last new entry. --------------------
$ws.Range("A19").Select() | Out-Null
